$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style from H1 (bold, bordered, centered) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Fill I2:I33 with constant 1
$ws.Range("I2:I33").Value = 1

# Fill J2:J33 with the same values as H2:H33 (copy values in one shot)
$ws.Range("H2:H33").Copy()
$ws.Range("J2:J33").PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = 0
